$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "in russia" -> "is_active"
$ws.Range("D1").Value = "is_active"

# Update the selected cell to D1
$ws.Range("D1").Select()
